$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.842.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.32%  '
$ws.Range("D3").Value = "'3.356.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.43%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'562.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.25%  '
$ws.Range("D6").Value = "'152.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.02%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = "'0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").Value = "'7.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("E10").Value = '  +3.74%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").Value = "'3.930.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.24%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("E14").Value = '  +2.67%  '
$ws.Range("E15").Value = '  +2.98%  '
$ws.Range("D16").Value = "'62.842.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.22%  '
$ws.Range("D17").Value = "'3.352.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.43%  '
$ws.Range("D18").Value = "'6.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("D19").Value = "'13.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.14%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = "'8.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = "'385.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = "'0.537"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.84%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  +5.21%  '
$ws.Range("D26").Value = "'8.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").Value = "'0.0₃0953"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.83%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("E29").Value = '  +5.76%  '
$ws.Range("E30").Value = '  +3.91%  '
$ws.Range("D31").Value = "'5.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.05%  '
$ws.Range("D32").Value = "'22.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.96%  '
$ws.Range("E33").Value = '  +6.85%  '
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("D35").Value = "'160.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.96%  '
$ws.Range("E36").Value = '  +7.90%  '
$ws.Range("E37").Value = '  +11.77%  '
$ws.Range("E38").Value = '  +3.50%  '
$ws.Range("E39").Value = '  +4.91%  '
$ws.Range("D40").Value = "'2.818.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("D41").Value = "'0.0311"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.38%  '
$ws.Range("D42").Value = "'0.747"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.47%  '
$ws.Range("D43").Value = "'40.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.15%  '
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("E45").Value = '  +3.04%  '
$ws.Range("D46").Value = "'3.397.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.35%  '
$ws.Range("D47").Value = "'21.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.27%  '
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("D50").Value = "'286.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.75%  '
$ws.Range("D51").Value = "'0.799"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.92%  '
